# Re-applies a re-scrape of betexplorer.com odds data:
#   - Several rows' match-data (columns F:V) got reshuffled among rows that
#     share the same match date/index (columns A:E are untouched).
#   - One brand-new match row (142) is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowVals($r) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$r").Value2
    }
    return $vals
}

function Set-RowVals($r, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $vals[$c]
    }
}

# --- snapshot every affected row's CURRENT (pre-edit) F:V content -----------
$rows = @(27,28,31, 70,71, 73,74,75, 89,90,91, 96,97, 103,104,107, 120,122, 132,133)
$before = @{}
foreach ($r in $rows) {
    $before[$r] = Get-RowVals $r
}

# --- target row -> source row (the row whose old content now lands here) ---
$mapping = @{
    27 = 28;  28 = 31;  31 = 27;
    70 = 71;  71 = 70;
    73 = 75;  74 = 73;  75 = 74;
    89 = 90;  90 = 91;  91 = 89;
    96 = 97;  97 = 96;
    103 = 104; 104 = 107; 107 = 103;
    120 = 122; 122 = 120;
    132 = 133; 133 = 132;
}

foreach ($r in $rows) {
    $srcRow = $mapping[$r]
    $srcVals = $before[$srcRow]
    Set-RowVals $r $srcVals
}

# --- append the new match as row 142 ----------------------------------------
# Duplicate row 141's formatting/styles down into row 142 first ...
$ws.Range("A141:V141").Copy($ws.Range("A142:V142"))

# ... then overwrite with the new match's own data.
$ws.Range("A142").Value2 = 141
$ws.Range("B142").Value2 = "poland"
$ws.Range("C142").Value2 = "iii-liga-group-iv"
$ws.Range("D142").Value2 = "2023-2024"
$ws.Range("E142").Value2 = 45258.5625
$ws.Range("F142").Value2 = "Orleta Radzyn"
$ws.Range("G142").Value2 = 1
$ws.Range("H142").Value2 = "Avia Swidnik"
$ws.Range("I142").Value2 = 3
$ws.Range("J142").Value2 = 4.4
$ws.Range("K142").Value2 = "25/11/2023 12:58"
$ws.Range("L142").Value2 = 4.4
$ws.Range("M142").Value2 = "25/11/2023 12:58"
$ws.Range("N142").Value2 = 4.04
$ws.Range("O142").Value2 = "25/11/2023 12:58"
$ws.Range("P142").Value2 = 4.04
$ws.Range("Q142").Value2 = "25/11/2023 12:58"
$ws.Range("R142").Value2 = 1.58
$ws.Range("S142").Value2 = "25/11/2023 12:58"
$ws.Range("T142").Value2 = 1.58
$ws.Range("U142").Value2 = "25/11/2023 12:58"
$ws.Range("V142").Value2 = "https://www.betexplorer.com/football/poland/iii-liga-group-iv/orleta-radzyn-avia-swidnik/QLOFQQtC/"
